$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.132.17"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "1.643.61"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'216.26"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'29.45"
$ws.Range("E8").Value = "  +9.49%  "
$ws.Range("E9").Value = "  +4.54%  "
$ws.Range("D10").Value = "'0.0616"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "'0.0917"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.876.80"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "1.641.14"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").Value = "'0.579"
$ws.Range("E14").Value = "  +7.29%  "
$ws.Range("D15").Value = "'9.62"
$ws.Range("E15").Value = "  +25.88%  "
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("D17").Value = "30.150.30"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "'65.01"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").Value = "'248.89"
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("D20").Value = "0.0₃0711"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'4.23"
$ws.Range("E22").Value = "  +6.09%  "
$ws.Range("D23").Value = "'10.01"
$ws.Range("E23").Value = "  +8.04%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").Value = "'159.73"
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("D26").Value = "'15.79"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("D28").Value = "'6.70"
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  +3.28%  "
$ws.Range("D31").Value = "'1.13"
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("D32").Value = "'3.42"
$ws.Range("E32").Value = "  +6.14%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "1.438.54"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'1.68"
$ws.Range("E35").Value = "  +8.51%  "
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0172"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.29"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'77.13"
$ws.Range("E40").Value = "  +16.76%  "
$ws.Range("D41").Value = "'0.563"
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.03"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.843"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "'55.53"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.0499"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D49").Value = "1.784.07"
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").Value = "'90.51"
$ws.Range("E50").Value = "  +5.01%  "
$ws.Range("D51").Value = "0.0₆0111"
$ws.Range("E51").Value = "  +6.13%  "
